$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ElementMap")
$ws2 = $wb.Worksheets.Item("TestCaseData")

# ---------------------------------------------------------------------------
# Sheet "ElementMap": insert a new "Header" section (3 rows) right after the
# LoginPage rows (new rows 5-7), moving the existing "lnk_goToCart" row into
# the new section and pushing everything else down accordingly.
# ---------------------------------------------------------------------------

# 1) Extend the row formatting down to the two brand new rows (26 & 27),
#    keeping the special "last row" formatting (row 25) on the new last
#    row (27), and the regular row formatting (row 24) on row 26.
$ws1.Range("A24:E24").Copy() | Out-Null
$ws1.Range("A26:E26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws1.Range("A25:E25").Copy() | Out-Null
$ws1.Range("A27:E27").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 2) Row 25 will no longer be the last row, so restore its formatting back to
#    the regular (non-last-row) style by copying from row 24.
$ws1.Range("A24:E24").Copy() | Out-Null
$ws1.Range("A25:E25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 3) Shift the data of rows 12-25 down to rows 14-27 (bottom-up so source
#    data is not overwritten before it is read).
for ($r = 25; $r -ge 12; $r--) {
    $v = $ws1.Range("A" + $r + ":E" + $r).Value2
    $ws1.Range("A" + ($r + 2) + ":E" + ($r + 2)).Value2 = $v
}

# 4) Shift the data of rows 5-10 down to rows 8-13.
for ($r = 10; $r -ge 5; $r--) {
    $v = $ws1.Range("A" + $r + ":E" + $r).Value2
    $ws1.Range("A" + ($r + 3) + ":E" + ($r + 3)).Value2 = $v
}

# 5) Populate the three new "Header" rows (5-7). Row 7 re-uses the
#    "lnk_goToCart" element that used to live under ProductsPage.
$ws1.Range("A5").Value = "Header"
$ws1.Range("B5").Value = "Button"
$ws1.Range("C5").Value = "btn_mainMenu"
$ws1.Range("D5").Value = "css"
$ws1.Range("E5").Value = "div#menu_button_container  > button"

$ws1.Range("A6").Value = "Header"
$ws1.Range("B6").Value = "Link"
$ws1.Range("C6").Value = "lnk_productsMenuLink"
$ws1.Range("D6").Value = "id"
$ws1.Range("E6").Value = "inventory_sidebar_link"

$ws1.Range("A7").Value = "Header"
$ws1.Range("B7").Value = "Link"
$ws1.Range("C7").Value = "lnk_goToCart"
$ws1.Range("D7").Value = "css"
$ws1.Range("E7").Value = "div#shopping_cart_container a[class*=shopping_cart_link]"

# ---------------------------------------------------------------------------
# Sheet "TestCaseData": the test case name changes casing.
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "buyAllTshirtsTest"

# ---------------------------------------------------------------------------
# View state: selections + active sheet change from ElementMap to
# TestCaseData.
# ---------------------------------------------------------------------------
$ws1.Range("E8").Select() | Out-Null
$ws2.Range("B7").Select() | Out-Null
$ws2.Activate() | Out-Null
